$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 20:05"

# Alemania (row 10) - updated counts
$ws.Range("B10").Value = 172905
$ws.Range("C10").Value = 329
$ws.Range("E10").Value = 18012
$ws.Range("G10").Value = 32
$ws.Range("H10").Value = 7693

# Canada (row 16) - updated counts
$ws.Range("B16").Value = 71099
$ws.Range("C16").Value = 1118
$ws.Range("D16").Value = 33707
$ws.Range("E16").Value = 32225
$ws.Range("G16").Value = 174
$ws.Range("H16").Value = 5167

# Israel (row 35) - updated counts
$ws.Range("B35").Value = 16529
$ws.Range("C35").Value = 23
$ws.Range("D35").Value = 12083
$ws.Range("E35").Value = 4186
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 260

# Egipto (row 49) - updated counts
$ws.Range("B49").Value = 10093
$ws.Range("C49").Value = 347
$ws.Range("D49").Value = 2326
$ws.Range("E49").Value = 7223
$ws.Range("G49").Value = 11
$ws.Range("H49").Value = 544

# Moldavia now overtakes Ghana (row 61) with fresh numbers,
# Ghana drops to row 62 keeping its previous numbers
$ws.Range("A61").Value = "Moldavia"
$ws.Range("B61").Value = 5154
$ws.Range("C61").Value = 159
$ws.Range("D61").Value = 2069
$ws.Range("E61").Value = 2903
$ws.Range("F61").Value = 251
$ws.Range("G61").Value = 7
$ws.Range("H61").Value = 182

$ws.Range("A62").Value = "Ghana"
$ws.Range("B62").Value = 5127
$ws.Range("C62").Value = 427
$ws.Range("D62").Value = 494
$ws.Range("E62").Value = 4611
$ws.Range("F62").Value = 5
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 22

# Somalia now overtakes Guatemala and Consejo Danes (row 93) with fresh numbers,
# Guatemala drops to row 94 and Consejo Danes drops to row 95, keeping previous numbers
$ws.Range("A93").Value = "Somalia"
$ws.Range("B93").Value = 1170
$ws.Range("C93").Value = 81
$ws.Range("D93").Value = 126
$ws.Range("E93").Value = 992
$ws.Range("F93").Value = 2
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 52

$ws.Range("A94").Value = "Guatemala"
$ws.Range("B94").Value = 1114
$ws.Range("C94").Value = 62
$ws.Range("D94").Value = 111
$ws.Range("E94").Value = 977
$ws.Range("F94").Value = 5
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 26

$ws.Range("A95").Value = "Consejo Danes para los Refugiados"
$ws.Range("B95").Value = 1102
$ws.Range("C95").Value = 78
$ws.Range("D95").Value = 146
$ws.Range("E95").Value = 912
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 44

# Georgia (row 119) - updated counts
$ws.Range("B119").Value = 642
$ws.Range("C119").Value = 4
$ws.Range("E119").Value = 282

# Mozambique (row 163) - updated counts
$ws.Range("B163").Value = 104
$ws.Range("C163").Value = 1
$ws.Range("E163").Value = 70
